$d = $word.ActiveDocument

$d.Content.Find.Execute("cddbb", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ccbdb", 2)
$d.Content.Find.Execute("acbaa", $true, $false, $false, $false, $false,
                         $true, 1, $false, "dcdcb", 2)
$d.Content.Find.Execute("babcb", $true, $false, $false, $false, $false,
                         $true, 1, $false, "dadda", 2)
$d.Content.Find.Execute("caadd", $true, $false, $false, $false, $false,
                         $true, 1, $false, "baddc", 2)
$d.Content.Find.Execute("dcaab", $true, $false, $false, $false, $false,
                         $true, 1, $false, "bdcad", 2)
